# Re-ran "resolve" and "classify+summarise" steps after changes to the
# mapping file. This zeroed-out / dropped several summary figures on the
# "Range Status" and "Species qualification" sheets, and re-derived the
# "High Priority break-up" sheet so it now only has a single ("IUCN") row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "Range Status" sheet: every species-count (col B) collapses to 0, and
# the accompanying percentage (col C) is dropped entirely (blank/cleared)
# for rows 2-7 (Historical, Very Restricted, Restricted, Moderate, Large,
# Very Large).
# ---------------------------------------------------------------------
$wsRange = $wb.Worksheets.Item("Range Status")
for ($r = 2; $r -le 7; $r++) {
    $wsRange.Cells.Item($r, 2).Value = 0
    $wsRange.Cells.Item($r, 3).ClearContents()
}

# ---------------------------------------------------------------------
# "Species qualification" sheet: Range Analysis (row 5) species count
# goes to 0.
# ---------------------------------------------------------------------
$wsQual = $wb.Worksheets.Item("Species qualification")
$wsQual.Range("B5").Value = 0

# ---------------------------------------------------------------------
# "High Priority break-up" sheet: the re-run collapsed the two rows
# (Range, IUCN) into one - the remaining row is labelled "IUCN" with
# updated figures, and the old row 3 is gone.
# ---------------------------------------------------------------------
$wsBreak = $wb.Worksheets.Item("High Priority break-up")
$wsBreak.Range("A2").Value = "IUCN"
$wsBreak.Range("B2").Value = 10
$wsBreak.Range("C2").Value = 100
$wsBreak.Range("D2").Value = 10
$wsBreak.Range("E2").Value = 100
$wsBreak.Rows.Item(3).Delete()
